# 2022 updating module order
# Swap module numbers between the two blocks of rows:
#   rows 39-48 (currently module 6) -> module 7
#   rows 49-53 (currently module 7) -> module 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 39; $r -le 48; $r++) {
    $ws.Cells.Item($r, 1).Value = 7
}

for ($r = 49; $r -le 53; $r++) {
    $ws.Cells.Item($r, 1).Value = 6
}
